$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("E120,E121,E122,E123,E125,E126,E127,E128,E129,E131,E132,E133,E134,E135,E137,E138,E139,E140,E141,E143,E144,E145,E146,E147,E149")
foreach ($c in $rng) {
    $c.NumberFormat = "@"
    $c.Value = "2018-04-28"
}

$rng = $ws.Range("E204,E205,E206,E207,E208,E210,E211,E212,E213,E214,E216,E217,E218,E219,E220,E222,E223,E224,E225,E226")
foreach ($c in $rng) {
    $c.NumberFormat = "@"
    $c.Value = "2018-06-14"
}

$rng = $ws.Range("E228,E229,E230,E231,E232,E234,E235,E236,E237,E238,E240,E241,E242,E243,E244,E246,E247,E248,E249,E250")
foreach ($c in $rng) {
    $c.NumberFormat = "@"
    $c.Value = "2019-06-22"
}

$rng = $ws.Range("E338,E350")
foreach ($c in $rng) {
    $c.NumberFormat = "@"
    $c.Value = "2015-04-10"
}

$rng = $ws.Range("E340,E352")
foreach ($c in $rng) {
    $c.NumberFormat = "@"
    $c.Value = "2015-05-23"
}

$rng = $ws.Range("E362,E374")
foreach ($c in $rng) {
    $c.NumberFormat = "@"
    $c.Value = "2016-04-14"
}

$rng = $ws.Range("E364,E376")
foreach ($c in $rng) {
    $c.NumberFormat = "@"
    $c.Value = "2016-05-21"
}

$rng = $ws.Range("E386,E398")
foreach ($c in $rng) {
    $c.NumberFormat = "@"
    $c.Value = "2017-04-14"
}

$rng = $ws.Range("E388,E400")
foreach ($c in $rng) {
    $c.NumberFormat = "@"
    $c.Value = "2017-05-24"
}

